$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was bumped
# from 46070 (2026-02-17) to 46072 (2026-02-19) for every data row (2-33).
for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46070) {
        $cell.Value2 = 46072
    }
}
